$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "92601306"
$ws.Range("D16").Value = "LUIS ENRIQUE BARRIOS SIERRA"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 25774
$ws.Range("G16").Value = 644350

$ws.Range("C17").Value = "17349569"
$ws.Range("D17").Value = "RUBIANO SANCHEZ MENDEZ"
$ws.Range("E17").Value = "1607"
$ws.Range("F17").Value = 25774
$ws.Range("G17").Value = 644350

$ws.Range("C18").Value = "1047446401"
$ws.Range("D18").Value = "FRANCISCO JAVIER RAMIREZ GAVIRIA"
$ws.Range("E18").Value = "1607"
$ws.Range("F18").Value = 25774
$ws.Range("G18").Value = 644350

$ws.Range("C19").Value = "73162890"
$ws.Range("D19").Value = "JUAN ALBERTO BARRIOS VASQUEZ"
$ws.Range("E19").Value = "1607"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

$ws.Range("C20").Value = "79397141"
$ws.Range("D20").Value = "OMAR RICARDO SEPULVEDA ADAMES"
$ws.Range("E20").Value = "1607"
$ws.Range("F20").Value = 25774
$ws.Range("G20").Value = 644350

$ws.Range("C21").Value = "9096062"
$ws.Range("D21").Value = "LUIS ALBERTO TORRES PUA"
$ws.Range("E21").Value = "1607"
$ws.Range("F21").Value = 25774
$ws.Range("G21").Value = 644350

$ws.Range("C22").Value = "1093140688"
$ws.Range("D22").Value = "MARLON RODRIGO GARCES CONTRERAS"
$ws.Range("E22").Value = "1607"
$ws.Range("F22").Value = 25774
$ws.Range("G22").Value = 644350

$ws.Range("C23").Value = "9167378"
$ws.Range("D23").Value = "MANUEL ARIAS BELLO"
$ws.Range("E23").Value = "1607"
$ws.Range("F23").Value = 25774
$ws.Range("G23").Value = 644350

$ws.Range("C24").Value = "92601306"
$ws.Range("D24").Value = "LUIS ENRIQUE BARRIOS SIERRA"
$ws.Range("E24").Value = "1608"
$ws.Range("F24").Value = 25774
$ws.Range("G24").Value = 644350

$ws.Range("C25").Value = "17349569"
$ws.Range("D25").Value = "RUBIANO SANCHEZ MENDEZ"
$ws.Range("E25").Value = "1608"
$ws.Range("F25").Value = 25774
$ws.Range("G25").Value = 644350

$ws.Range("C26").Value = "1047446401"
$ws.Range("D26").Value = "FRANCISCO JAVIER RAMIREZ GAVIRIA"
$ws.Range("E26").Value = "1608"
$ws.Range("F26").Value = 25774
$ws.Range("G26").Value = 644350

$ws.Range("C27").Value = "73162890"
$ws.Range("D27").Value = "JUAN ALBERTO BARRIOS VASQUEZ"
$ws.Range("E27").Value = "1608"
$ws.Range("F27").Value = 40000
$ws.Range("G27").Value = 1000000

$ws.Range("C28").Value = "79397141"
$ws.Range("D28").Value = "OMAR RICARDO SEPULVEDA ADAMES"
$ws.Range("E28").Value = "1608"
$ws.Range("F28").Value = 25774
$ws.Range("G28").Value = 644350

$ws.Range("C29").Value = "9096062"
$ws.Range("D29").Value = "LUIS ALBERTO TORRES PUA"
$ws.Range("E29").Value = "1608"
$ws.Range("F29").Value = 25774
$ws.Range("G29").Value = 644350

$ws.Range("C30").Value = "1093140688"
$ws.Range("D30").Value = "MARLON RODRIGO GARCES CONTRERAS"
$ws.Range("E30").Value = "1608"
$ws.Range("F30").Value = 25774
$ws.Range("G30").Value = 644350

$ws.Range("C31").Value = "9167378"
$ws.Range("D31").Value = "MANUEL ARIAS BELLO"
$ws.Range("E31").Value = "1608"
$ws.Range("F31").Value = 25774
$ws.Range("G31").Value = 644350

